$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 181 (shifts old rows 181-195 down to 182-196)
$ws.Range("A181").EntireRow.Insert()

# Populate the newly inserted row 181 with the new record
$ws.Range("A181").Value = 8
$ws.Range("B181").Value = "Terminal La Palmera de La Serena"
$ws.Range("C181").Value = "Coquimbo"
$ws.Range("D181").Value = 44615
$ws.Range("E181").Value = 4
$ws.Range("F181").Value = 100112021
$ws.Range("G181").Value = "Ají"
$ws.Range("H181").Value = "Americana (o)"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 560
$ws.Range("K181").Value = 12000
$ws.Range("L181").Value = 13000
$ws.Range("M181").Value = 12500
$ws.Range("N181").Value = "$/caja 15 kilos"
$ws.Range("O181").Value = "Provincia de Limarí"
$ws.Range("P181").Value = 833
$ws.Range("Q181").Value = 15
$ws.Range("R181").Value = "Hortaliza"
